$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.851.58"
$ws.Range("E2").Value = "  -0.55%  "

$ws.Range("D3").Value = "3.955.12"
$ws.Range("E3").Value = "  -2.44%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'607.33"
$ws.Range("E5").Value = "  +1.59%  "

$ws.Range("D6").Value = "'171.21"
$ws.Range("E6").Value = "  +11.26%  "

$ws.Range("E7").Value = "  -1.23%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").Value = "'0.786"
$ws.Range("E9").Value = "  +3.47%  "

$ws.Range("D10").Value = "'0.180"
$ws.Range("E10").Value = "  +5.89%  "

$ws.Range("D11").Value = "'56.30"
$ws.Range("E11").Value = "  +4.95%  "

$ws.Range("E12").Value = "  +0.49%  "

$ws.Range("D13").Value = "'11.51"
$ws.Range("E13").Value = "  +4.09%  "

$ws.Range("D14").Value = "4.581.83"
$ws.Range("E14").Value = "  -2.63%  "

$ws.Range("D15").Value = "3.961.65"
$ws.Range("E15").Value = "  -2.33%  "

$ws.Range("D16").Value = "'21.24"
$ws.Range("E16").Value = "  +2.15%  "

$ws.Range("E17").Value = "  -1.75%  "

$ws.Range("E18").Value = "  -0.79%  "

$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").Value = "'0.131"
$ws.Range("E19").Value = "  -0.45%  "

$ws.Range("B20").Value = "WrappedBTC"
$ws.Range("C20").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D20").Value = "72.732.33"
$ws.Range("E20").Value = "  -0.72%  "

$ws.Range("D21").Value = "'444.16"
$ws.Range("E21").Value = "  +0.42%  "

$ws.Range("D22").Value = "'4.87"
$ws.Range("E22").Value = "  +0.99%  "

$ws.Range("D23").Value = "'95.81"
$ws.Range("E23").Value = "  -1.97%  "

$ws.Range("E24").Value = "  -5.04%  "

$ws.Range("E25").Value = "  -1.76%  "

$ws.Range("E26").Value = "  -0.91%  "

$ws.Range("D27").Value = "'11.40"
$ws.Range("E27").Value = "  -0.06%  "

$ws.Range("E28").Value = "  -1.38%  "

$ws.Range("E29").Value = "  -4.33%  "

$ws.Range("D30").Value = "'35.81"
$ws.Range("E30").Value = "  -3.27%  "

$ws.Range("D31").Value = "'7.93"
$ws.Range("E31").Value = "  +0.05%  "

$ws.Range("D32").Value = "'13.87"
$ws.Range("E32").Value = "  +1.50%  "

$ws.Range("D33").Value = "'50.02"
$ws.Range("E33").Value = "  +2.31%  "

$ws.Range("D34").Value = "'0.128"
$ws.Range("E34").Value = "  -4.17%  "

$ws.Range("D35").Value = "0.0₃0998"
$ws.Range("E35").Value = "  +13.89%  "

$ws.Range("D36").Value = "'69.33"
$ws.Range("E36").Value = "  -1.58%  "

$ws.Range("D37").Value = "'636.97"
$ws.Range("E37").Value = "  -7.61%  "

$ws.Range("D38").Value = "'0.430"
$ws.Range("E38").Value = "  -4.11%  "

$ws.Range("D39").Value = "'3.45"
$ws.Range("E39").Value = "  +2.52%  "

$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.146"
$ws.Range("E40").Value = "  -1.34%  "

$ws.Range("B41").Value = "Dai"
$ws.Range("C41").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  -0.08%  "

$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.06%  "

$ws.Range("E43").Value = "  -3.41%  "

$ws.Range("D44").Value = "'10.58"
$ws.Range("E44").Value = "  -6.56%  "

$ws.Range("D45").Value = "'3.16"
$ws.Range("E45").Value = "  +41.57%  "

$ws.Range("E46").Value = "  -2.14%  "

$ws.Range("E47").Value = "  -2.30%  "

$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "'3.40"
$ws.Range("E48").Value = "  -0.14%  "

$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").Value = "'2.88"
$ws.Range("E49").Value = "  -14.45%  "

$ws.Range("D50").Value = "'0.000285"
$ws.Range("E50").Value = "  +5.47%  "

$ws.Range("D51").Value = "2.839.13"
$ws.Range("E51").Value = "  +1.36%  "
